$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -12.54149999999999
$ws.Range("A9").Value = -21.99880000000001
$ws.Range("C9").Value = -13.4183
$ws.Range("D9").Value = -8.917199999999999
$ws.Range("C11").Value = -12.78449999999999
$ws.Range("A13").Value = -22.3125
$ws.Range("A16").Value = -21.5417
$ws.Range("C16").Value = -13.3833
$ws.Range("A18").Value = -22.1492
$ws.Range("A20").Value = -21.70829999999998
$ws.Range("D22").Value = -7.841200000000002
$ws.Range("C23").Value = -12.3123
$ws.Range("C24").Value = -12.8608
$ws.Range("A26").Value = -21.96459999999998
$ws.Range("C26").Value = -13.1154
$ws.Range("A27").Value = -21.8748
$ws.Range("D27").Value = -8.652199999999999
$ws.Range("A29").Value = -20.70749999999997
$ws.Range("D29").Value = -7.214999999999995
$ws.Range("D32").Value = -6.644599999999993
$ws.Range("C34").Value = -11.81300000000001
$ws.Range("A35").Value = -21.89059999999998
$ws.Range("C35").Value = -11.94499999999999
$ws.Range("A36").Value = -21.88189999999999
$ws.Range("D37").Value = -7.801000000000001
$ws.Range("D38").Value = -7.835100000000006
$ws.Range("D39").Value = -6.841899999999994
$ws.Range("D41").Value = -7.839600000000002
$ws.Range("C44").Value = -12.1635
$ws.Range("A45").Value = -21.42149999999998
$ws.Range("D45").Value = -7.60109999999999
$ws.Range("C48").Value = -12.8857
$ws.Range("D48").Value = -8.590699999999998
$ws.Range("C49").Value = -14.0488
$ws.Range("D51").Value = -8.760199999999999
$ws.Range("C52").Value = -10.8961
$ws.Range("A55").Value = -22.3628
$ws.Range("D56").Value = -9.224800000000004
$ws.Range("A57").Value = -22.39180000000001
$ws.Range("D57").Value = -7.921199999999994
$ws.Range("D61").Value = -7.993099999999997
$ws.Range("D64").Value = -7.103799999999995
$ws.Range("C66").Value = -11.019
$ws.Range("C67").Value = -10.5139
$ws.Range("A69").Value = -21.5816
$ws.Range("C73").Value = -10.68410000000001
$ws.Range("D75").Value = -8.756299999999996
$ws.Range("A76").Value = -19.61999999999998
$ws.Range("A78").Value = -19.86309999999999
$ws.Range("C78").Value = -10.83050000000001
$ws.Range("C80").Value = -11.20140000000001
$ws.Range("A82").Value = -22.07710000000001
$ws.Range("D82").Value = -8.781200000000005
$ws.Range("A83").Value = -21.56109999999999
$ws.Range("D90").Value = -7.342099999999995
$ws.Range("C91").Value = -13.21009999999999
$ws.Range("A93").Value = -21.38970000000001
$ws.Range("D93").Value = -6.90309999999999
$ws.Range("A97").Value = -21.56020000000001
$ws.Range("C97").Value = -10.41900000000001
$ws.Range("C99").Value = -11.6729
$ws.Range("D102").Value = -7.954900000000007
$ws.Range("C104").Value = -11.41810000000001
$ws.Range("D105").Value = -8.080099999999998